# Add "sprint 4" technologies to the radar (Kong API Gateway, Guava)
# as two new rows at the bottom of Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 58: Kong API Gateway -------------------------------------------
# Set F58 ("sprint4") first, then A58, then E58 so that the new shared
# strings are appended to sharedStrings.xml in the same order as the
# reference edit (sprint4, Kong API Gateway, description, ...).
$ws.Range("F58").Value = "sprint4"
$ws.Range("A58").Value = "Kong API Gateway"
$ws.Range("B58").Value = "Evaluer"
$ws.Range("C58").Value = "Plateforme"
# Copy the existing "TRUE" text cell so column D keeps being stored as the
# literal shared string "TRUE" (and not auto-converted to a boolean).
$ws.Range("D57").Copy($ws.Range("D58"))
$ws.Range("E58").Value = "API Gateway Open Souce et haute performance. Basé sur nGinx. Implémentation minimale mais essentielles d'une Gateway API."

# --- Row 59: Guava --------------------------------------------------------
$ws.Range("A59").Value = "Guava"
$ws.Range("B59").Value = "Adopeter"
$ws.Range("C59").Value = "Outils et Librairies"
$ws.Range("D57").Copy($ws.Range("D59"))
$ws.Range("E59").Value = "Librairie fourni par Google proposant toute une panoplie d'outils facilitant la vie du dévelopeur (traitement des collections, validations des paramètres, etc..)"
$ws.Range("F59").Value = "sprint4"

# Style the two new rows with the new green font used for sprint4 entries.
# (The cells start out with the default Calibri/12 font, so only the color
# needs to change to reproduce the new font exactly.)
$ws.Range("A58:F59").Font.Color = 5287936

# Update the visible selection to match the saved view.
$ws.Range("E63").Select() | Out-Null
